$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "World" series (column F) is being dropped from the table; deleting
# the entire column shifts Rest of Europe & Asia / MENA / South & SE Asia /
# Sub-Saharan Africa (old G:J) one slot to the left (new F:I) along with
# their header labels, matching the trimmed shared-strings table.
$ws.Columns.Item(6).Delete()

# Leave the selection where the author's cursor ended up after the edit.
$ws.Range("AA6").Select()
